$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $text) {
    $cell = $table.Cell($row, $col)
    $cellRange = $cell.Range
    $cellRange.End = $cellRange.End - 1
    $cellRange.Text = $text
}

# Row 2: Banana
Set-CellText $table 2 2 '$ 0.73'
Set-CellText $table 2 3 '1'
Set-CellText $table 2 4 '$ 0.73'

# Row 3: Strawberry
Set-CellText $table 3 2 '$ 0.28'
Set-CellText $table 3 3 '8'
Set-CellText $table 3 4 '$ 2.24'

# Row 4: Chicken
Set-CellText $table 4 2 '$ 0.99'
Set-CellText $table 4 4 '$ 5.94'

# Row 5: Bread
Set-CellText $table 5 2 '$ 0.16'
Set-CellText $table 5 3 '2'
Set-CellText $table 5 4 '$ 0.32'

# Row 6: Eggs
Set-CellText $table 6 2 '$ 0.27'
Set-CellText $table 6 3 '3'
Set-CellText $table 6 4 '$ 0.81'

# Row 7: Salad
Set-CellText $table 7 2 '$ 0.41'
Set-CellText $table 7 3 '2'
Set-CellText $table 7 4 '$ 0.82'
